# Update "想去人数" (number of people interested) values for several
# events across the "展览" (Exhibition), "演出" (Performance) and
# "全部类型" (All types) sheets, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1799
$ws1.Range("F6").Value  = 671
$ws1.Range("F7").Value  = 41
$ws1.Range("F13").Value = 165
$ws1.Range("F17").Value = 111
$ws1.Range("F18").Value = 5096
$ws1.Range("F21").Value = 116
$ws1.Range("F22").Value = 2272
$ws1.Range("F25").Value = 2121

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 84

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1799
$ws4.Range("F6").Value  = 671
$ws4.Range("F7").Value  = 41
$ws4.Range("F13").Value = 165
$ws4.Range("F17").Value = 111
$ws4.Range("F18").Value = 5096
$ws4.Range("F19").Value = 84
$ws4.Range("F23").Value = 116
$ws4.Range("F24").Value = 2272
$ws4.Range("F27").Value = 2121
